$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rationale text (column B) and lag values (columns C/D) ---

# Row 2: beninvert_cpue
$ws.Range("B2").Value = "Represents prey items across ontogeny. Missing polychetes, which are preferential juvenile crab prey. Limited prey may influence survival to recruitment the following year "
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# Row 3: cp_extent
$ws.Range("C3").Value = 4

# Row 4: Mean_AO
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 3

# Row 6: bcs_imm
$ws.Range("B6").Value = "Highest prevelances are in small juveniles- affecting survival to recruitment  "
$ws.Range("C6").Value = 3

# Row 8: temp_occ_imm
$ws.Range("B8").Value = "Juvenile snow crab habitat; high temperatures may indirectly or directly affect survival to recruitment "
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 1

# Row 9: JanFeb_ice
$ws.Range("B9").Value = "Los of sea ice likely decreases both food quality and quantity supplied to the benthos for early juveniles- Copeman et. al research"
$ws.Range("D9").Value = 1

# Row 11: Pcod_consumption
$ws.Range("D11").Value = 1

# Row 12: Chla
$ws.Range("B12").Value = "Diatoms are a key prey item for larval stages"
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 3

# Row 13: juv_condition
$ws.Range("B13").Value = "Males collected for this indicator are larger than 50-65mm, assuming similar mechanism between cohorts for reduced condition within a year. Poor condition suggests poor survival to recruitment the following year "
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1

# Row 16: model_recruit methodology note
$ws.Range("B16").Value = "Recruitment output from 2022 approved model. Note that the final few years of recruitment estimates are fairly unreliable. Size crab recruiting to the model are 25-40mm carapace width, and therefore likely 3-4 years post-settlement"

# --- Column B width adjustment ---
# Target stored width is 84.109375 characters; the COM ColumnWidth setter
# quantizes internally to the nearest 1/6 character (pixel rounding), so we
# pick the input that yields the closest achievable stored value (84.1667).
$ws.Columns.Item(2).ColumnWidth = 83.33

# --- Selected cell ---
$ws.Range("C6").Select()
